$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths (B, C, D widened to fit the new user-story text) ---
# Target "characters" widths from the authored workbook are 20.140625 / 29 /
# 93.7109375. This runtime's ColumnWidth setter stores width = round(input*6)/6
# + 5/6 (points internally, /6 back to "characters" on export) -- a 1/6-char
# grid -- so feeding the literal target is NOT generally the closest hit
# (e.g. literal 29 lands on 29.833, not 29). Instead we feed the input value
# that snaps onto the grid point nearest each target: exact for column C
# (29 is reachable), within ~0.03-0.04 chars for B/D (not on the grid, so
# unreachable exactly by any input through this API).
$ws.Columns.Item(2).ColumnWidth = 19.333333333333332
$ws.Columns.Item(3).ColumnWidth = 28.166666666666668
$ws.Columns.Item(4).ColumnWidth = 92.83333333333333

# --- New user-story rows (A2:G6) ---
# Columns: ID, As a/an, I want to, so that, notes, priority, status
$stories = @(
    [pscustomobject]@{ Id = 1; AsA = "user"; IWantTo = "Enter Student Numbers using a textbox"; SoThat = "The student numbers enter the database"; Notes = $null;                                      Priority = "High";   Status = "Done" },
    [pscustomobject]@{ Id = 2; AsA = "user"; IWantTo = "View the time-ins of today"; SoThat = "I view the time-ins of that day"; Notes = $null;                                                       Priority = "Normal"; Status = "In-progress" },
    [pscustomobject]@{ Id = 3; AsA = "user"; IWantTo = "View the current student numbers timed-in"; SoThat = "So that I can keep track of the current time-ins"; Notes = $null;                       Priority = "High";   Status = "In-progress" },
    [pscustomobject]@{ Id = 4; AsA = "user"; IWantTo = "To be able to use a barcode scanner"; SoThat = "My job will be easier"; Notes = "Simulate using Android barcode Scanner";                     Priority = "High";   Status = "Done" },
    [pscustomobject]@{ Id = 5; AsA = "user"; IWantTo = "To be able to generate report from data"; SoThat = "I can see the trends in the data provided"; Notes = $null;                                 Priority = "Normal"; Status = "In-progress" }
)

$row = 2
foreach ($story in $stories) {
    $ws.Cells.Item($row, 1).Value = $story.Id
    $ws.Cells.Item($row, 2).Value = $story.AsA
    $ws.Cells.Item($row, 3).Value = $story.IWantTo
    $ws.Cells.Item($row, 4).Value = $story.SoThat
    if ($story.Notes) {
        $ws.Cells.Item($row, 5).Value = $story.Notes
    }
    $ws.Cells.Item($row, 6).Value = $story.Priority
    $ws.Cells.Item($row, 7).Value = $story.Status
    $row++
}

# Selection ends up on the first empty row below the table, same as Excel
# would leave it after typing the last row and pressing Enter.
$ws.Range("A7").Select()
